# fix bug template import
# Appends a bold red "(*)" marker to the four header labels that must be
# filled in when importing this template, fixes the active selection and
# page orientation left over from the author's last save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-HeaderWithStar {
    param(
        [string]$cellRef,
        [string]$prefix,
        [string]$starText
    )

    $full = $prefix + $starText
    $ws.Range($cellRef).Value = $full

    # Characters() is 1-indexed; bold+red only the "(*)" run.
    $startPos = $prefix.Length + 1
    $len = $starText.Length
    $chars = $ws.Range($cellRef).Characters($startPos, $len)
    $chars.Font.Bold = $true
    $chars.Font.Color = 255
}

Set-HeaderWithStar "A1" "Mã phiếu " "(*)"
Set-HeaderWithStar "J1" "Line / STT dòng " "(*)"
Set-HeaderWithStar "K1" "Line / Sản phẩm" " (*)"
Set-HeaderWithStar "L1" "Line / Số lượng xuất " "(*)"

# Widen the Line columns now that the headers carry the extra "(*)" text.
$ws.Columns.Item(10).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(11).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(12).EntireColumn.AutoFit() | Out-Null

# Restore portrait page setup.
$ws.PageSetup.Orientation = 1

# Leave selection where the author left it on save.
$ws.Range("G10").Select() | Out-Null
